$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue 2 4 '30.772.60'
Set-TextValue 2 5 '  +1.70%  '
Set-TextValue 3 4 '2.108.81'
Set-TextValue 3 5 '  +5.57%  '
Set-TextValue 4 4 '1.001'
Set-TextValue 4 5 '  -0.13%  '
Set-TextValue 5 4 '331.51'
Set-TextValue 5 5 '  +2.53%  '
Set-TextValue 6 4 '1.001'
Set-TextValue 6 5 '  -0.06%  '
Set-TextValue 7 4 '0.5299'
Set-TextValue 7 5 '  +3.92%  '
Set-TextValue 8 4 '0.4363'
Set-TextValue 8 5 '  +5.93%  '
Set-TextValue 9 4 '0.08957'
Set-TextValue 9 5 '  +3.03%  '
Set-TextValue 10 4 '47.19'
Set-TextValue 10 5 '  +10.43%  '
Set-TextValue 11 5 '  +3.13%  '
Set-TextValue 12 4 '24.76'
Set-TextValue 12 5 '  +0.20%  '
Set-TextValue 13 4 '2.104.43'
Set-TextValue 13 5 '  +5.35%  '
Set-TextValue 14 4 '6.721'
Set-TextValue 14 5 '  +3.18%  '
Set-TextValue 15 4 '7.772'
Set-TextValue 15 5 '  +4.78%  '
Set-TextValue 16 4 '96.74'
Set-TextValue 16 5 '  +2.90%  '
Set-TextValue 17 5 '  -0.04%  '
Set-TextValue 18 5 '  +1.45%  '
Set-TextValue 19 4 '0.06672'
Set-TextValue 19 5 '  +1.94%  '
Set-TextValue 20 4 '19.03'
Set-TextValue 20 5 '  +0.73%  '
Set-TextValue 21 4 '1.001'
Set-TextValue 21 5 '  -0.03%  '
Set-TextValue 22 4 '6.297'
Set-TextValue 22 5 '  +2.64%  '
Set-TextValue 23 4 '30.830.85'
Set-TextValue 23 5 '  +1.71%  '
Set-TextValue 24 5 '  +3.37%  '
Set-TextValue 25 4 '2.349.36'
Set-TextValue 25 5 '  +5.36%  '
Set-TextValue 26 5 '  +3.06%  '
Set-TextValue 27 4 '22.58'
Set-TextValue 28 4 '2.581'
Set-TextValue 28 5 '  +7.35%  '
Set-TextValue 29 4 '162.13'
Set-TextValue 29 5 '  -0.73%  '
Set-TextValue 30 4 '132.96'
Set-TextValue 30 5 '  +1.21%  '
Set-TextValue 31 4 '1.194'
Set-TextValue 31 5 '  +5.01%  '
Set-TextValue 32 4 '0.1079'
Set-TextValue 32 5 '  +2.73%  '
Set-TextValue 33 4 '6.155'
Set-TextValue 33 5 '  +1.80%  '
Set-TextValue 34 4 '1.550'
Set-TextValue 34 5 '  +17.03%  '
Set-TextValue 35 4 '3.922'
Set-TextValue 35 5 '  +2.39%  '
Set-TextValue 36 4 '0.02596'
Set-TextValue 36 5 '  +3.83%  '
Set-TextValue 37 4 '9.644'
Set-TextValue 37 5 '  +7.75%  '
Set-TextValue 38 4 '5.533'
Set-TextValue 38 5 '  +3.34%  '
Set-TextValue 39 4 '0.06741'
Set-TextValue 39 5 '  +2.46%  '
Set-TextValue 40 5 '  +3.57%  '
Set-TextValue 41 4 '0.2274'
Set-TextValue 41 5 '  +3.52%  '
Set-TextValue 42 4 '0.6835'
Set-TextValue 43 4 '1.242'
Set-TextValue 43 5 '  +1.29%  '
Set-TextValue 44 4 '1.000'
Set-TextValue 44 5 '  -0.12%  '
Set-TextValue 45 2 'Decentraland'
Set-TextValue 45 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 45 4 '0.6410'
Set-TextValue 45 5 '  +4.37%  '
Set-TextValue 46 2 'EnergySwap'
Set-TextValue 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 46 4 '13.92'
Set-TextValue 46 5 '  +1.82%  '
Set-TextValue 47 5 '  +0.92%  '
Set-TextValue 48 4 '3.646'
Set-TextValue 48 5 '  -0.41%  '
Set-TextValue 49 4 '1.258'
Set-TextValue 49 5 '  -0.49%  '
Set-TextValue 50 4 '82.85'
Set-TextValue 50 5 '  +3.73%  '
Set-TextValue 51 4 '1.191'
Set-TextValue 51 5 '  +7.94%  '
